# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on the last data row
# (row 14) of the per-locale sheets, reflecting a newer report-generation
# run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D14").Value = "2016-03-08 06:25:04"
$wsZhCn.Range("G14").Value = "2016-03-08 06:25:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D14").Value = "2016-03-08 06:25:07"
$wsDeDe.Range("G14").Value = "2016-03-08 06:25:26"
